$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 1077.1052
$ws.Range("I4").Value = 1025.8889
$ws.Range("K4").Value = 1025.8889
$ws.Range("M4").Value = -911.8888999999999
# Row 19
$ws.Range("H19").Value = 1291.1666
$ws.Range("I19").Value = 1441.6666
$ws.Range("J19").Value = 1215.9166
$ws.Range("K19").Value = 1441.6666
$ws.Range("L19").Value = 1215.9166
$ws.Range("M19").Value = -1266.6666
$ws.Range("N19").Value = -1565.9166
# Row 33
$ws.Range("H33").Value = 33054.742
$ws.Range("I33").Value = 40213.6
$ws.Range("J33").Value = 3226.1667
$ws.Range("K33").Value = 40213.6
$ws.Range("L33").Value = 3226.1667
$ws.Range("M33").Value = -39984.6
$ws.Range("N33").Value = -3684.1667
# Row 51
$ws.Range("H51").Value = 12246.182
$ws.Range("I51").Value = 22040.2
$ws.Range("K51").Value = 22040.2
$ws.Range("M51").Value = -21556.2
# Row 98
$ws.Range("H98").Value = 1255.3636
$ws.Range("I98").Value = 1356.5555
$ws.Range("K98").Value = 1356.5555
$ws.Range("M98").Value = 141.4445000000001
# Row 113
$ws.Range("H113").Value = 113011.664
$ws.Range("I113").Value = 113011.664
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 113011.664
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -109757.664
$ws.Range("N113").ClearContents()
# Row 122
$ws.Range("H122").Value = 1255.3636
$ws.Range("I122").Value = 1356.5555
$ws.Range("K122").Value = 4069.6665
$ws.Range("M122").Value = -1619.6665
# Row 129
$ws.Range("H129").Value = 874.48486
$ws.Range("I129").Value = 565.875
$ws.Range("J129").Value = 973.24
$ws.Range("K129").Value = 1697.625
$ws.Range("L129").Value = 2919.72
$ws.Range("M129").Value = 3302.375
$ws.Range("N129").Value = -12919.72

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3024.6365
$ws.Range("I45").Value = 2531.6667
$ws.Range("K45").Value = 2531.6667
$ws.Range("M45").Value = -2154.6667
# Row 122
$ws.Range("H122").Value = 5915.8184
$ws.Range("I122").Value = 7413.3335
$ws.Range("J122").Value = 4118.8
$ws.Range("K122").Value = 22240.0005
$ws.Range("L122").Value = 12356.4
$ws.Range("M122").Value = -19790.0005
$ws.Range("N122").Value = -17256.4
# Row 139
$ws.Range("H139").Value = 54500.555
$ws.Range("J139").Value = 54500.555
$ws.Range("L139").Value = 54500.555
$ws.Range("N139").Value = -64780.555

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 13
$ws.Range("H13").Value = 42000
$ws.Range("J13").Value = 42000
$ws.Range("L13").Value = 42000
$ws.Range("N13").Value = -42336
# Row 20
$ws.Range("H20").Value = 26458.512
$ws.Range("I20").Value = 41149.92
$ws.Range("J20").Value = 3503.1875
$ws.Range("K20").Value = 41149.92
$ws.Range("L20").Value = 3503.1875
$ws.Range("M20").Value = -40902.92
$ws.Range("N20").Value = -3997.1875

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 13526
$ws.Range("I31").Value = 19467.814
$ws.Range("K31").Value = 19467.814
$ws.Range("M31").Value = -19172.814
# Row 34
$ws.Range("H34").Value = 13526
$ws.Range("I34").Value = 19467.814
$ws.Range("K34").Value = 19467.814
$ws.Range("M34").Value = -19265.814
# Row 58
$ws.Range("H58").Value = 6170.18
$ws.Range("I58").Value = 978.2432
$ws.Range("J58").Value = 20947.23
$ws.Range("K58").Value = 978.2432
$ws.Range("L58").Value = 20947.23
$ws.Range("M58").Value = -775.2432
$ws.Range("N58").Value = -21353.23
# Row 132
$ws.Range("H132").Value = 3061
$ws.Range("I132").Value = 2747.7856
$ws.Range("J132").Value = 3687.4285
$ws.Range("K132").Value = 8243.356800000001
$ws.Range("L132").Value = 11062.2855
$ws.Range("M132").Value = -5713.356800000001
$ws.Range("N132").Value = -16122.2855
# Row 136
$ws.Range("H136").Value = 6170.18
$ws.Range("I136").Value = 978.2432
$ws.Range("J136").Value = 20947.23
$ws.Range("K136").Value = 2934.7296
$ws.Range("L136").Value = 62841.69
$ws.Range("M136").Value = -384.7296000000001
$ws.Range("N136").Value = -67941.69
# Row 140
$ws.Range("H140").Value = 49800
$ws.Range("J140").Value = 49800
$ws.Range("L140").Value = 49800
$ws.Range("N140").Value = -60160

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 655.9474
$ws.Range("I113").Value = 591.25
$ws.Range("J113").Value = 766.8570999999999
$ws.Range("K113").Value = 1773.75
$ws.Range("L113").Value = 2300.5713
$ws.Range("M113").Value = 396.25
$ws.Range("N113").Value = -6640.5713
# Row 122
$ws.Range("H122").Value = 13537.25
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
# Row 131
$ws.Range("H131").Value = 808.58
$ws.Range("J131").Value = 876.4773
$ws.Range("L131").Value = 2629.4319
$ws.Range("N131").Value = -12709.4319
# Row 132
$ws.Range("H132").Value = 3524.0476
$ws.Range("I132").Value = 3437.5
$ws.Range("J132").Value = 3577.3076
$ws.Range("K132").Value = 30937.5
$ws.Range("L132").Value = 32195.7684
$ws.Range("M132").Value = -28407.5
$ws.Range("N132").Value = -37255.7684
# Row 141
$ws.Range("H141").Value = 3503.1785
$ws.Range("J141").Value = 3595.4092
$ws.Range("L141").Value = 10786.2276
$ws.Range("N141").Value = -21146.2276

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 49974.16
$ws.Range("I70").Value = 73172.92999999999
$ws.Range("K70").Value = 73172.92999999999
$ws.Range("M70").Value = -72902.92999999999
# Row 73
$ws.Range("H73").Value = 49974.16
$ws.Range("I73").Value = 73172.92999999999
$ws.Range("K73").Value = 73172.92999999999
$ws.Range("M73").Value = -72236.92999999999
# Row 132
$ws.Range("H132").Value = 3272.7144
$ws.Range("I132").Value = 3359
$ws.Range("J132").Value = 3057
$ws.Range("K132").Value = 10077
$ws.Range("L132").Value = 9171
$ws.Range("M132").Value = -7547
$ws.Range("N132").Value = -14231

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1693.5834
$ws.Range("I93").Value = 1717
$ws.Range("J93").Value = 1670.1666
$ws.Range("K93").Value = 1717
$ws.Range("L93").Value = 1670.1666
$ws.Range("M93").Value = -469
$ws.Range("N93").Value = -4166.1666
# Row 132
$ws.Range("H132").Value = 7042.6875
$ws.Range("I132").Value = 7791.273
$ws.Range("J132").Value = 5395.8
$ws.Range("K132").Value = 23373.819
$ws.Range("L132").Value = 16187.4
$ws.Range("M132").Value = -20843.819
$ws.Range("N132").Value = -21247.4

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 125
$ws.Range("H125").Value = 31702.6
$ws.Range("J125").Value = 31702.6
$ws.Range("L125").Value = 31702.6
$ws.Range("N125").Value = -41542.6

